# Append/update the acquisition-timestamp column ("取得日時") on the first
# sheet ("ランサーズ") for all existing data rows (rows 2-10), changing the
# previous run's timestamp "2025-10-26 18:22:53" to the new run's timestamp
# "2025-10-26 18:32:10", as described in the commit message:
#   "Append: 2025-10-26 18:32 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2:A10").Value = "2025-10-26 18:32:10"
